$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (POX/C)
$ws.Range("B2").Value = 100
$ws.Range("C2").Value = 10
$ws.Range("D2").Value = 119.5996398454005
$ws.Range("E2").Value = 80.40036015459947

# Row 3 (C/A)
$ws.Range("B3").Value = 0.1
$ws.Range("C3").Value = 0.01
$ws.Range("D3").Value = 0.1195996398454005
$ws.Range("E3").Value = 0.08040036015459946

# Row 4 (POX/M)
$ws.Range("B4").Value = 0.001
$ws.Range("C4").Value = 0.0001
$ws.Range("D4").Value = 0.001195996398454005
$ws.Range("E4").Value = 0.0008040036015459946
